$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.978.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.748.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5205"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.80%  "

$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.756.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07024"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6458"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.539"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.989.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006635"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.980.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.166"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.670"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.158"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.507"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.841"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08302"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.672"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.445"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04474"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9884"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6135"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.681"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.936"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3874"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.088"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7355"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05472"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.327"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.644"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
